$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '64.450.23'
$ws.Cells.Item(2, 5).Value = '  -3.68%  '
$ws.Cells.Item(3, 4).Value = '3.027.03'
$ws.Cells.Item(3, 5).Value = '  -6.10%  '
$ws.Cells.Item(4, 5).Value = '  +0.36%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '556.47'
$ws.Cells.Item(5, 5).Value = '  -6.28%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '139.04'
$ws.Cells.Item(6, 5).Value = '  -8.60%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.998'
$ws.Cells.Item(7, 5).Value = '  -0.10%  '
$ws.Cells.Item(8, 4).Value = '3.016.00'
$ws.Cells.Item(8, 5).Value = '  -6.18%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.478'
$ws.Cells.Item(9, 5).Value = '  -12.59%  '
$ws.Cells.Item(10, 5).Value = '  -12.08%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '6.01'
$ws.Cells.Item(11, 5).Value = '  -9.61%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.445'
$ws.Cells.Item(12, 5).Value = '  -11.55%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '33.89'
$ws.Cells.Item(13, 5).Value = '  -13.69%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.0000211'
$ws.Cells.Item(14, 5).Value = '  -14.34%  '
$ws.Cells.Item(15, 4).Value = '3.528.52'
$ws.Cells.Item(15, 5).Value = '  -5.72%  '
$ws.Cells.Item(16, 4).Value = '64.572.64'
$ws.Cells.Item(16, 5).Value = '  -3.53%  '
$ws.Cells.Item(17, 5).Value = '  -3.87%  '
$ws.Cells.Item(18, 4).Value = '3.044.70'
$ws.Cells.Item(18, 5).Value = '  -5.51%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '480.84'
$ws.Cells.Item(19, 5).Value = '  -9.85%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '6.40'
$ws.Cells.Item(20, 5).Value = '  -11.42%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '13.31'
$ws.Cells.Item(21, 5).Value = '  -11.74%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.652'
$ws.Cells.Item(22, 5).Value = '  -14.66%  '
$ws.Cells.Item(23, 5).Value = '  -14.16%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '12.41'
$ws.Cells.Item(24, 5).Value = '  -10.51%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '77.59'
$ws.Cells.Item(25, 5).Value = '  -9.85%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.999'
$ws.Cells.Item(26, 5).Value = '  -0.17%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '2.69'
$ws.Cells.Item(27, 5).Value = '  -15.56%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.05'
$ws.Cells.Item(28, 5).Value = '  -7.46%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '7.45'
$ws.Cells.Item(29, 5).Value = '  -9.07%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '25.51'
$ws.Cells.Item(30, 5).Value = '  -13.46%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '2.55'
$ws.Cells.Item(31, 5).Value = '  -3.58%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.12'
$ws.Cells.Item(32, 5).Value = '  -0.90%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '522.99'
$ws.Cells.Item(33, 5).Value = '  -3.68%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.01'
$ws.Cells.Item(34, 5).Value = '  +0.42%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '5.22'
$ws.Cells.Item(35, 5).Value = '  -9.36%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '51.69'
$ws.Cells.Item(36, 5).Value = '  -2.77%  '
$ws.Cells.Item(37, 5).Value = '  -12.95%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.0404'
$ws.Cells.Item(38, 5).Value = '  -4.62%  '
$ws.Cells.Item(39, 2).Value = 'Hedera'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.0781'
$ws.Cells.Item(39, 5).Value = '  -10.64%  '
$ws.Cells.Item(40, 2).Value = 'Kaspa'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.119'
$ws.Cells.Item(40, 5).Value = '  -6.61%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '8.06'
$ws.Cells.Item(41, 5).Value = '  -14.11%  '
$ws.Cells.Item(42, 4).Value = '2.801.97'
$ws.Cells.Item(42, 5).Value = '  -4.91%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '2.39'
$ws.Cells.Item(43, 5).Value = '  -9.66%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.998'
$ws.Cells.Item(44, 5).Value = '  -0.19%  '
$ws.Cells.Item(45, 5).Value = '  -12.18%  '
$ws.Cells.Item(46, 5).Value = '  -8.22%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '115.69'
$ws.Cells.Item(47, 5).Value = '  -5.48%  '
$ws.Cells.Item(48, 4).Value = '0.0₃0504'
$ws.Cells.Item(48, 5).Value = '  -14.57%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.105'
$ws.Cells.Item(49, 5).Value = '  -8.75%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '23.28'
$ws.Cells.Item(50, 5).Value = '  -12.30%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '1.98'
$ws.Cells.Item(51, 5).Value = '  -18.22%  '
